$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the floating-point rounding of the timestamp in A12 (same moment, different
# representation of the fractional seconds: 19:00:21.472215 -> 19:00:21.472000)
$ws.Range("A12").Value2 = 45863.79191518518

# Append the new scheduled-task reading as row 13
$ws.Range("A13").NumberFormat = $ws.Range("A12").NumberFormat
$ws.Range("A13").Value2 = 45863.83357264067
$ws.Range("B13").Value2 = 2025
$ws.Range("C13").Value2 = 30
$ws.Range("D13").Value2 = 14
$ws.Range("E13").Value2 = 85.37
$ws.Range("F13").Value2 = 0
$ws.Range("G13").Value2 = 3.51
$ws.Range("H13").Value = "E"
$ws.Range("I13").Value2 = 0
$ws.Range("J13").Value = "20:00:20"
